$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in missing "Problem Number" values for existing rows ---
$ws.Range("B55").Value = 34
$ws.Range("B56").Value = 36
$ws.Range("B57").Value = 46
$ws.Range("B58").Value = 48
$ws.Range("B59").Value = 49
$ws.Range("B60").Value = 54

# --- Add the new row 63 entry for LeetCode 189 - Rotate Array ---
$ws.Range("A63").Value = "Array"
$ws.Range("B63").Value = 189
$ws.Range("C63").Value = "189-Rotate Array"
$ws.Range("D63").Value = "Medium"
$ws.Range("E63").Value = "Cyclic replacement"
$ws.Range("F63").Value = "Failed"
$ws.Range("G63").Value = "O(n) time, O(1) memory"
$ws.Range("H63").Value = "Reverse Array"
$ws.Range("I63").Value = "O(N) time, O(1) memory"
$ws.Range("J63").Value = "no"
$ws.Range("K63").Value = "no"
$ws.Range("L63").Value = "I spent a lot of time trying to come up with a cyclic approach which was really tricky to account for the edge cases.`nThe trick to this problem was to relaize you can reverse the array, and then reverse the array 2 more times in 2 separate portions to obtain the end result"
$ws.Range("M63").Value = "1 hour 30 minutes"
$ws.Range("N63").Value = "yes"

$ws.Rows.Item(63).RowHeight = 225

# --- Update the view: scroll position & active selection ---
$excel.ActiveWindow.ScrollRow = 60
$ws.Range("A64").Select()
